# Apply "more work towards final product" changes:
#  - add carrier word (column D) to the four practice rows
#  - tag the first two generic pairs as unique_video and the next two as unique_audio (column J)
#  - add eight new rows (9-16) describing the unique_video / unique_audio carrier pairing
#    (columns C = kind, D = carrier)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Practice rows (p1-p4): add carrier word in column D
$ws.Range("D2").Value = "can"
$ws.Range("D3").Value = "where"
$ws.Range("D4").Value = "do"
$ws.Range("D5").Value = "look"

# Generic pair rows 1-4 (sheet rows 6-9): mark pair_kind in column J
$ws.Range("J6").Value = "unique_video"
$ws.Range("J7").Value = "unique_video"
$ws.Range("J8").Value = "unique_audio"
$ws.Range("J9").Value = "unique_audio"

# New rows for numbers 9-16 (sheet rows 14-21): kind (C) and carrier (D)
$ws.Range("C14").Value = "unique_video"
$ws.Range("D14").Value = "can"

$ws.Range("C15").Value = "unique_video"
$ws.Range("D15").Value = "can"

$ws.Range("C16").Value = "unique_video"
$ws.Range("D16").Value = "do"

$ws.Range("C17").Value = "unique_video"
$ws.Range("D17").Value = "do"

$ws.Range("C18").Value = "unique_audio"
$ws.Range("D18").Value = "look"

$ws.Range("C19").Value = "unique_audio"
$ws.Range("D19").Value = "look"

$ws.Range("C20").Value = "unique_audio"
$ws.Range("D20").Value = "where"

$ws.Range("C21").Value = "unique_audio"
$ws.Range("D21").Value = "where"
